$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2541899441340782
$ws.Range("C2").Value = 0.4664804469273743
$ws.Range("J2").Value = 0.0223463687150838
$ws.Range("P2").Value = 0.1955307262569832
$ws.Range("S2").Value = 0.06145251396648044
$ws.Range("B3").Value = 0.01169590643274854
$ws.Range("C3").Value = 0.01169590643274854
$ws.Range("J3").Value = 0.01754385964912281
$ws.Range("P3").Value = 0.7192982456140351
$ws.Range("S3").Value = 0.239766081871345
$ws.Range("J4").Value = 0.05555555555555555
$ws.Range("P4").Value = 0.6851851851851852
$ws.Range("S4").Value = 0.2592592592592592
$ws.Range("B6").Value = 0.04102564102564103
$ws.Range("D6").Value = 0.02051282051282051
$ws.Range("F6").Value = 0.09230769230769231
$ws.Range("J6").Value = 0.2615384615384616
$ws.Range("O6").Value = 0.005128205128205128
$ws.Range("Q6").Value = 0.1282051282051282
$ws.Range("R6").Value = 0.05128205128205128
$ws.Range("S6").Value = 0.4
$ws.Range("B7").Value = 0.1352657004830918
$ws.Range("D7").Value = 0.02898550724637681
$ws.Range("F7").Value = 0.06763285024154589
$ws.Range("J7").Value = 0.1400966183574879
$ws.Range("O7").Value = 0.01932367149758454
$ws.Range("Q7").Value = 0.1739130434782609
$ws.Range("R7").Value = 0.06763285024154589
$ws.Range("S7").Value = 0.3671497584541063
$ws.Range("B8").Value = 0.103448275862069
$ws.Range("D8").Value = 0.01379310344827586
$ws.Range("F8").Value = 0.03908045977011494
$ws.Range("J8").Value = 0.1402298850574713
$ws.Range("O8").Value = 0.02068965517241379
$ws.Range("Q8").Value = 0.1494252873563219
$ws.Range("R8").Value = 0.09425287356321839
$ws.Range("S8").Value = 0.439080459770115
$ws.Range("B9").Value = 0.08125
$ws.Range("D9").Value = 0.025
$ws.Range("F9").Value = 0.04375
$ws.Range("J9").Value = 0.1125
$ws.Range("O9").Value = 0.03125
$ws.Range("Q9").Value = 0.19375
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.4125
$ws.Range("B10").Value = 0.1372854914196568
$ws.Range("D10").Value = 0.02730109204368175
$ws.Range("E10").Value = 0.0015600624024961
$ws.Range("F10").Value = 0.05226209048361934
$ws.Range("J10").Value = 0.1193447737909516
$ws.Range("O10").Value = 0.0187207488299532
$ws.Range("Q10").Value = 0.2223088923556942
$ws.Range("R10").Value = 0.08736349453978159
$ws.Range("S10").Value = 0.3338533541341654
$ws.Range("G11").Value = 0.1501597444089457
$ws.Range("J11").Value = 0.09584664536741214
$ws.Range("K11").Value = 0.1789137380191693
$ws.Range("L11").Value = 0.5623003194888179
$ws.Range("S11").Value = 0.01277955271565495
$ws.Range("G12").Value = 0.774869109947644
$ws.Range("J12").Value = 0.162303664921466
$ws.Range("K12").Value = 0.01047120418848168
$ws.Range("L12").Value = 0.01570680628272251
$ws.Range("S12").Value = 0.03664921465968586
$ws.Range("G13").Value = 0.7428571428571429
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.05714285714285714
$ws.Range("F15").Value = 0.01401869158878505
$ws.Range("H15").Value = 0.1355140186915888
$ws.Range("I15").Value = 0.0514018691588785
$ws.Range("J15").Value = 0.3925233644859813
$ws.Range("K15").Value = 0.06542056074766354
$ws.Range("M15").Value = 0.01401869158878505
$ws.Range("O15").Value = 0.06542056074766354
$ws.Range("S15").Value = 0.2616822429906542
$ws.Range("F16").Value = 0.03083700440528634
$ws.Range("H16").Value = 0.1409691629955947
$ws.Range("I16").Value = 0.08370044052863436
$ws.Range("J16").Value = 0.3876651982378855
$ws.Range("K16").Value = 0.13215859030837
$ws.Range("M16").Value = 0.00881057268722467
$ws.Range("O16").Value = 0.07048458149779736
$ws.Range("S16").Value = 0.145374449339207
$ws.Range("F17").Value = 0.02528735632183908
$ws.Range("H17").Value = 0.1494252873563219
$ws.Range("I17").Value = 0.0735632183908046
$ws.Range("J17").Value = 0.4528735632183908
$ws.Range("K17").Value = 0.08735632183908046
$ws.Range("M17").Value = 0.01379310344827586
$ws.Range("O17").Value = 0.07586206896551724
$ws.Range("S17").Value = 0.1218390804597701
$ws.Range("F18").Value = 0.02083333333333333
$ws.Range("I18").Value = 0.07291666666666667
$ws.Range("J18").Value = 0.4375
$ws.Range("K18").Value = 0.09375
$ws.Range("M18").Value = 0.015625
$ws.Range("N18").Value = 0.005208333333333333
$ws.Range("O18").Value = 0.0625
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.01845906902086678
$ws.Range("H19").Value = 0.2207062600321027
$ws.Range("I19").Value = 0.06581059390048154
$ws.Range("J19").Value = 0.3635634028892456
$ws.Range("K19").Value = 0.1171749598715891
$ws.Range("M19").Value = 0.01605136436597111
$ws.Range("O19").Value = 0.05617977528089887
$ws.Range("S19").Value = 0.1420545746388443
